$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update boardinghouse entry
$ws.Range("B2").Value = "Kantunan sa UEP"
$ws.Range("C2").Value = "Ivan Joseph G. Arang"
# leading apostrophe forces text storage so the leading zero in the
# contact number is preserved; reset the style afterwards so no
# extraneous number-format style is left applied to the cell
$ws.Range("D2").Value = "'09385050074"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "Avelino Street"
$ws.Range("F2").Value = "Zone 3"

# Row 3: update boardinghouse entry (Street/Zone reuse same values as row 2)
$ws.Range("B3").Value = "House of Lanister"
$ws.Range("C3").Value = "Melvin Dionisio"
$ws.Range("D3").Value = "'09518015683"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "Avelino Street"
$ws.Range("F3").Value = "Zone 3"
